# Update "想去人数" (number of people interested) counts on several rows
# across the 展览, 演出 and 全部类型 sheets, per the source-site refresh
# (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1015
$ws1.Range("F4").Value  = 13384
$ws1.Range("F5").Value  = 39
$ws1.Range("F6").Value  = 1010
$ws1.Range("F8").Value  = 1723
$ws1.Range("F11").Value = 69
$ws1.Range("F14").Value = 13365
$ws1.Range("F16").Value = 586
$ws1.Range("F17").Value = 8904
$ws1.Range("F18").Value = 2
$ws1.Range("F19").Value = 7978
$ws1.Range("F20").Value = 243
$ws1.Range("F26").Value = 19
$ws1.Range("F27").Value = 1015
$ws1.Range("F28").Value = 12
$ws1.Range("F30").Value = 391
$ws1.Range("F31").Value = 199
$ws1.Range("F32").Value = 157

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 32

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1015
$ws4.Range("F5").Value  = 13384
$ws4.Range("F6").Value  = 39
$ws4.Range("F7").Value  = 1010
$ws4.Range("F9").Value  = 1723
$ws4.Range("F12").Value = 69
$ws4.Range("F15").Value = 13365
$ws4.Range("F17").Value = 586
$ws4.Range("F18").Value = 8904
$ws4.Range("F19").Value = 2
$ws4.Range("F20").Value = 7978
$ws4.Range("F21").Value = 243
$ws4.Range("F27").Value = 19
$ws4.Range("F28").Value = 1015
$ws4.Range("F29").Value = 12
$ws4.Range("F31").Value = 32
$ws4.Range("F33").Value = 391
$ws4.Range("F34").Value = 199
$ws4.Range("F35").Value = 157
